$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename grid entries (version 1d variable renaming)
$ws.Range("A3").Value2 = "create_synthetic_voltage_control_lv_network"
$ws.Range("A4").Value2 = "example_multivoltage"
$ws.Range("A5").Value2 = "example_simple"
$ws.Range("A6").Value2 = "kb_extrem_dorfnetz"
$ws.Range("A7").Value2 = "mv_oberrhein"
$ws.Range("B7").Value2 = "OPF Not Converged"

# Remove the now-obsolete trailing rows (8-25)
$ws.Rows("8:25").Delete()
